# Update "想去人数" (Wanted-to-go count) column F on all sheets to reflect
# newly generated numbers for the gh-pages output.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 622
$ws.Range("F5").Value = 2745
$ws.Range("F7").Value = 213
$ws.Range("F10").Value = 6206
$ws.Range("F11").Value = 4
$ws.Range("F13").Value = 76
$ws.Range("F14").Value = 4979
$ws.Range("F16").Value = 534
$ws.Range("F17").Value = 2591
$ws.Range("F18").Value = 1339
$ws.Range("F19").Value = 1509
$ws.Range("F20").Value = 1213
$ws.Range("F21").Value = 298
$ws.Range("F22").Value = 116
$ws.Range("F24").Value = 1046
$ws.Range("F26").Value = 386
$ws.Range("F27").Value = 522
$ws.Range("F28").Value = 1356
$ws.Range("F29").Value = 1019
$ws.Range("F30").Value = 2083
$ws.Range("F32").Value = 572
$ws.Range("F33").Value = 16
$ws.Range("F35").Value = 243
$ws.Range("F36").Value = 1486
$ws.Range("F38").Value = 1037
$ws.Range("F39").Value = 112
$ws.Range("F40").Value = 548
$ws.Range("F41").Value = 11
$ws.Range("F43").Value = 2245
$ws.Range("F44").Value = 2535
$ws.Range("F46").Value = 119
$ws.Range("F47").Value = 269
$ws.Range("F48").Value = 102
$ws.Range("F49").Value = 84

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 17
$ws.Range("F8").Value = 316
$ws.Range("F10").Value = 84
$ws.Range("F15").Value = 149
$ws.Range("F22").Value = 323
$ws.Range("F23").Value = 352
$ws.Range("F29").Value = 17
$ws.Range("F32").Value = 4

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1693
$ws.Range("F8").Value = 1461
$ws.Range("F9").Value = 1809
$ws.Range("F10").Value = 2481
$ws.Range("F11").Value = 826

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1693
$ws.Range("F6").Value = 622
$ws.Range("F8").Value = 2745
$ws.Range("F9").Value = 213
$ws.Range("F10").Value = 1461
$ws.Range("F12").Value = 2481
$ws.Range("F13").Value = 6206
$ws.Range("F14").Value = 826
$ws.Range("F17").Value = 4979
$ws.Range("F18").Value = 2591
$ws.Range("F19").Value = 1509
$ws.Range("F20").Value = 1213
$ws.Range("F21").Value = 298
$ws.Range("F22").Value = 116
$ws.Range("F24").Value = 1046
$ws.Range("F26").Value = 84
$ws.Range("F27").Value = 386
$ws.Range("F28").Value = 1356
$ws.Range("F29").Value = 1019
$ws.Range("F30").Value = 2083
$ws.Range("F32").Value = 572
$ws.Range("F33").Value = 243
$ws.Range("F35").Value = 1486
$ws.Range("F36").Value = 1037
$ws.Range("F37").Value = 548
$ws.Range("F39").Value = 323
$ws.Range("F42").Value = 2245
$ws.Range("F43").Value = 2535
$ws.Range("F44").Value = 119
$ws.Range("F45").Value = 269
$ws.Range("F46").Value = 102
$ws.Range("F47").Value = 84
$ws.Range("F48").Value = 4
